# Applies the edits described by the commit "Datamanager and Finetuning fitting"
# to the "Thoughts on Plotting window" document.
#
# Edits are applied from the bottom of the document upwards so that the
# (fixed, original) 1-based Paragraphs(...) indices used below stay valid
# right up until the very last step, which deletes a paragraph near the
# top and shifts everything after it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new, empty paragraph right after the "Thinking of making data
#    manager ..." paragraph (original paragraph 17).
# ---------------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
$p17.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 2. Fix the typo "tings" -> "things" (inside paragraph 17 only) and split
#    the run around "things" the same way the saved document does:
#    "...and other " | "things" | " about fits..."
# ---------------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
$rng = $p17.Range.Duplicate()
$found = $rng.Find.Execute("tings", $true, $false, $false, $false, $false, $true, 1, $false, "things", 2)

$p17 = $d.Paragraphs(17)
$rng2 = $p17.Range.Duplicate()
$found2 = $rng2.Find.Execute("things", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Bold = 1
    $rng2.Bold = 0
}

# ---------------------------------------------------------------------------
# 3. "Plotting -window : ... remove lines from plots." / " Option should be
#    given ... being by the side ..." -> re-split (paragraph 15 only) so the
#    second run starts at "the side and cannot be changed." instead.
# ---------------------------------------------------------------------------
$p15 = $d.Paragraphs(15)
$fullText = "Plotting -window : should have option to remove lines from plots. Option should be given in data manager, I think this should be separate windows or actual tabs with the plotting window being by the side and cannot be changed. But for now I think they should be windows that spawn when you want to see them with the plotting window being original GUI. They should save state so can be closed without changing anything"
$rng = $p15.Range.Duplicate()
$found = $rng.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, $fullText, 2)

$p15 = $d.Paragraphs(15)
$target = "the side and cannot be changed. But for now I think they should be windows that spawn when you want to see them with the plotting window being original GUI. They should save state so can be closed without changing anything"
$rng2 = $p15.Range.Duplicate()
$found2 = $rng2.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Bold = 1
    $rng2.Bold = 0
}

# ---------------------------------------------------------------------------
# 4. "Custom plot should essential just give ..." -> "... essentially just
#    give ..." (paragraph 14 only) and split the run the same way the saved
#    document does: "Custom plot should essential" | "ly" | " just give..."
# ---------------------------------------------------------------------------
$p14 = $d.Paragraphs(14)
$rng = $p14.Range.Duplicate()
$found = $rng.Find.Execute("essential just", $true, $false, $false, $false, $false, $true, 1, $false, "essentially just", 2)

$p14 = $d.Paragraphs(14)
$rng2 = $p14.Range.Duplicate()
$found2 = $rng2.Find.Execute("ly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Bold = 1
    $rng2.Bold = 0
}

# ---------------------------------------------------------------------------
# 5. Remove the first "Details of fitting tab" bullet ("Create your custom
#    function", paragraph 8). The remaining bullets shift up, which also
#    realises the text changes the diff shows for the next three bullets
#    (their text is unchanged; only their position moves).
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
if ($p8.Range.Text -like "Create your custom function*") {
    $p8.Range.Delete()
}
